# Week 13 stat logging update
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# M.Jones (row 2): RZATT 13 -> 15
$rushing.Cells.Item(2, 5).Value = 15

# D.Harris (row 3): 2DATT 96 -> 102, 3DATT 47 -> 50, RZATT 10 -> 12
$rushing.Cells.Item(3, 3).Value = 102
$rushing.Cells.Item(3, 4).Value = 50
$rushing.Cells.Item(3, 5).Value = 12

# R.Stevenson (row 5): 2DATT 48 -> 59, 3DATT 25 -> 35, RZATT 2 -> 5, last col 17 -> 20
$rushing.Cells.Item(5, 3).Value = 59
$rushing.Cells.Item(5, 4).Value = 35
$rushing.Cells.Item(5, 5).Value = 5
$rushing.Cells.Item(5, 6).Value = 20

# B.Bolden (row 7): RZATT 11 -> 15, last col 4 -> 5
$rushing.Cells.Item(7, 5).Value = 15
$rushing.Cells.Item(7, 6).Value = 5

# N.Agholor (row 8): 3DATT 0 -> 1
$rushing.Cells.Item(8, 4).Value = 1

# K.Bourne (row 9): 2DATT 4 -> 5
$rushing.Cells.Item(9, 3).Value = 5

# J.Smith (row 11): 2DATT 1 -> 2
$rushing.Cells.Item(11, 3).Value = 2

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# B.Bolden (row 5): Short Target 27 -> 28, Short Comp 25 -> 26
$receiving.Cells.Item(5, 3).Value = 28
$receiving.Cells.Item(5, 4).Value = 26

# N.Agholor (row 6): Short Target 34 -> 35
$receiving.Cells.Item(6, 3).Value = 35

# J.Smith (row 12): Short Target 32 -> 33, Short Comp 21 -> 22
$receiving.Cells.Item(12, 3).Value = 33
$receiving.Cells.Item(12, 4).Value = 22
